# Auto-generated: apply numeric corrections to Leve profit calc columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 914.0278
$ws.Range("I15").Value = 914.0278
$ws.Range("K15").Value = 2742.0834
$ws.Range("M15").Value = -2573.0834
$ws.Range("H62").Value = 4998
$ws.Range("I62").Value = 4998
$ws.Range("K62").Value = 4998
$ws.Range("M62").Value = -4374
$ws.Range("H65").Value = 4998
$ws.Range("I65").Value = 4998
$ws.Range("K65").Value = 24990
$ws.Range("M65").Value = -21870
$ws.Range("H98").Value = 3127.15
$ws.Range("I98").Value = 2366.1538
$ws.Range("K98").Value = 2366.1538
$ws.Range("M98").Value = -868.1538
$ws.Range("H112").Value = 3186.0488
$ws.Range("I112").Value = 2132.3333
$ws.Range("J112").Value = 3269.2368
$ws.Range("K112").Value = 6396.999899999999
$ws.Range("L112").Value = 9807.7104
$ws.Range("M112").Value = -5288.999899999999
$ws.Range("N112").Value = -12023.7104
$ws.Range("H116").Value = 5997.2856
$ws.Range("I116").Value = 6129.6665
$ws.Range("K116").Value = 6129.6665
$ws.Range("M116").Value = -2687.6665
$ws.Range("H122").Value = 3127.15
$ws.Range("I122").Value = 2366.1538
$ws.Range("K122").Value = 7098.4614
$ws.Range("M122").Value = -4648.4614
$ws.Range("H127").Value = 1297.5714
$ws.Range("I127").Value = 1297.5714
$ws.Range("K127").Value = 3892.7142
$ws.Range("M127").Value = 1067.2858
$ws.Range("H130").Value = 116222.5
$ws.Range("J130").Value = 116222.5
$ws.Range("L130").Value = 116222.5
$ws.Range("N130").Value = -126262.5
$ws.Range("H132").Value = 1943.28
$ws.Range("I132").Value = 1742.0952
$ws.Range("K132").Value = 5226.2856
$ws.Range("M132").Value = -2696.2856
$ws.Range("H135").Value = 1978.2
$ws.Range("I135").Value = 1788.1111
$ws.Range("J135").Value = 2263.3333
$ws.Range("K135").Value = 16092.9999
$ws.Range("L135").Value = 20369.9997
$ws.Range("M135").Value = -13557.9999
$ws.Range("N135").Value = -25439.9997
$ws.Range("H137").Value = 1375958.8
$ws.Range("I137").Value = 4330.9453
$ws.Range("J137").Value = 5567044
$ws.Range("K137").Value = 12992.8359
$ws.Range("L137").Value = 16701132
$ws.Range("M137").Value = -10442.8359
$ws.Range("N137").Value = -16706232
$ws.Range("H138").Value = 5616.1
$ws.Range("I138").Value = 7926.533
$ws.Range("J138").Value = 4229.84
$ws.Range("K138").Value = 23779.599
$ws.Range("L138").Value = 12689.52
$ws.Range("M138").Value = -18639.599
$ws.Range("N138").Value = -22969.52

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 190399.47
$ws.Range("I32").Value = 200643.14
$ws.Range("J32").Value = 19671.666
$ws.Range("K32").Value = 200643.14
$ws.Range("L32").Value = 19671.666
$ws.Range("M32").Value = -200356.14
$ws.Range("N32").Value = -20245.666
$ws.Range("H61").Value = 1002887.5
$ws.Range("I61").Value = 27241.762
$ws.Range("J61").Value = 3994867.8
$ws.Range("K61").Value = 27241.762
$ws.Range("L61").Value = 3994867.8
$ws.Range("M61").Value = -27029.762
$ws.Range("N61").Value = -3995291.8
$ws.Range("H74").Value = 397322.47
$ws.Range("I74").Value = 3506.0195
$ws.Range("J74").Value = 2907902.5
$ws.Range("K74").Value = 3506.0195
$ws.Range("L74").Value = 2907902.5
$ws.Range("M74").Value = -2632.0195
$ws.Range("N74").Value = -2909650.5
$ws.Range("H77").Value = 397322.47
$ws.Range("I77").Value = 3506.0195
$ws.Range("J77").Value = 2907902.5
$ws.Range("K77").Value = 17530.0975
$ws.Range("L77").Value = 14539512.5
$ws.Range("M77").Value = -13162.0975
$ws.Range("N77").Value = -14548248.5
$ws.Range("H132").Value = 3124.2222
$ws.Range("I132").Value = 1402.625
$ws.Range("J132").Value = 4501.5
$ws.Range("K132").Value = 4207.875
$ws.Range("L132").Value = 13504.5
$ws.Range("M132").Value = -1677.875
$ws.Range("N132").Value = -18564.5
$ws.Range("H136").Value = 1002887.5
$ws.Range("I136").Value = 27241.762
$ws.Range("J136").Value = 3994867.8
$ws.Range("K136").Value = 81725.28599999999
$ws.Range("L136").Value = 11984603.4
$ws.Range("M136").Value = -79175.28599999999
$ws.Range("N136").Value = -11989703.4
$ws.Range("H137").Value = 82298
$ws.Range("J137").Value = 82298
$ws.Range("L137").Value = 82298
$ws.Range("N137").Value = -92498

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 26061.6
$ws.Range("I99").Value = 26061.6
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 26061.6
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -24563.6
$ws.Range("N99").ClearContents()
$ws.Range("H134").Value = 75002990
$ws.Range("I134").Value = 3198.8572
$ws.Range("J134").Value = 180002700
$ws.Range("K134").Value = 9596.571599999999
$ws.Range("L134").Value = 540008100
$ws.Range("M134").Value = -7061.571599999999
$ws.Range("N134").Value = -540013170

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 8000
$ws.Range("I17").Value = 8000
$ws.Range("K17").Value = 8000
$ws.Range("M17").Value = -7826
$ws.Range("H38").Value = 8450
$ws.Range("I38").Value = 5000
$ws.Range("K38").Value = 5000
$ws.Range("M38").Value = -4623
$ws.Range("H39").Value = 5035.4287
$ws.Range("I39").Value = 2349.6
$ws.Range("K39").Value = 2349.6
$ws.Range("M39").Value = -1958.6
$ws.Range("H46").Value = 8450
$ws.Range("I46").Value = 5000
$ws.Range("K46").Value = 5000
$ws.Range("M46").Value = -4789
$ws.Range("H49").Value = 5035.4287
$ws.Range("I49").Value = 2349.6
$ws.Range("K49").Value = 2349.6
$ws.Range("M49").Value = -2167.6
$ws.Range("H58").Value = 2288.5833
$ws.Range("I58").Value = 2752
$ws.Range("K58").Value = 2752
$ws.Range("M58").Value = -2549
$ws.Range("H132").Value = 2948.5908
$ws.Range("I132").Value = 2731.3333
$ws.Range("J132").Value = 3099
$ws.Range("K132").Value = 8193.999899999999
$ws.Range("L132").Value = 9297
$ws.Range("M132").Value = -5663.999899999999
$ws.Range("N132").Value = -14357
$ws.Range("H136").Value = 2288.5833
$ws.Range("I136").Value = 2752
$ws.Range("K136").Value = 8256
$ws.Range("M136").Value = -5706

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1491.0667
$ws.Range("J132").Value = 2149.8333
$ws.Range("L132").Value = 19348.4997
$ws.Range("N132").Value = -24408.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1117809.9
$ws.Range("I11").Value = 1117809.9
$ws.Range("K11").Value = 1117809.9
$ws.Range("M11").Value = -1117670.9
$ws.Range("H19").Value = 4998
$ws.Range("I19").Value = 4996
$ws.Range("K19").Value = 4996
$ws.Range("M19").Value = -4708
$ws.Range("H97").Value = 77446.96000000001
$ws.Range("I97").Value = 59219.59
$ws.Range("J97").Value = 111876.445
$ws.Range("K97").Value = 59219.59
$ws.Range("L97").Value = 111876.445
$ws.Range("M97").Value = -58723.59
$ws.Range("N97").Value = -112868.445
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830
$ws.Range("H132").Value = 803012.0600000001
$ws.Range("I132").Value = 1208.8572
$ws.Range("J132").Value = 1027516.94
$ws.Range("K132").Value = 3626.5716
$ws.Range("L132").Value = 3082550.82
$ws.Range("M132").Value = -1096.5716
$ws.Range("N132").Value = -3087610.82

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 12259.5
$ws.Range("I46").Value = 15268.75
$ws.Range("K46").Value = 15268.75
$ws.Range("M46").Value = -15080.75
$ws.Range("H61").Value = 1457.5834
$ws.Range("J61").Value = 1361.25
$ws.Range("L61").Value = 1361.25
$ws.Range("N61").Value = -1765.25
$ws.Range("H113").Value = 1457.5834
$ws.Range("J113").Value = 1361.25
$ws.Range("L113").Value = 1361.25
$ws.Range("N113").Value = -5701.25
$ws.Range("H122").Value = 3548.077
$ws.Range("I122").Value = 2357.5
$ws.Range("K122").Value = 7072.5
$ws.Range("M122").Value = -4622.5
$ws.Range("H130").Value = 77539
$ws.Range("J130").Value = 77539
$ws.Range("L130").Value = 77539
$ws.Range("N130").Value = -87579
$ws.Range("H132").Value = 2853.8928
$ws.Range("I132").Value = 2765.762
$ws.Range("J132").Value = 3118.2856
$ws.Range("K132").Value = 8297.286
$ws.Range("L132").Value = 9354.856800000001
$ws.Range("M132").Value = -5767.286
$ws.Range("N132").Value = -14414.8568
$ws.Range("H134").Value = 69985
$ws.Range("J134").Value = 69985
$ws.Range("L134").Value = 69985
$ws.Range("N134").Value = -80125
$ws.Range("H136").Value = 33352.363
$ws.Range("I136").Value = 64875.875
$ws.Range("J136").Value = 3683.1765
$ws.Range("K136").Value = 194627.625
$ws.Range("L136").Value = 11049.5295
$ws.Range("M136").Value = -192077.625
$ws.Range("N136").Value = -16149.5295
$ws.Range("H138").Value = 34999
$ws.Range("J138").Value = 34999
$ws.Range("L138").Value = 34999
$ws.Range("N138").Value = -45279
$ws.Range("H141").Value = 139990
$ws.Range("J141").Value = 139990
$ws.Range("L141").Value = 139990
$ws.Range("N141").Value = -150350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 10000
$ws.Range("I37").Value = 10000
$ws.Range("K37").Value = 10000
$ws.Range("M37").Value = -9797
$ws.Range("H46").Value = 68139
$ws.Range("J46").Value = 68139
$ws.Range("L46").Value = 68139
$ws.Range("N46").Value = -68601
$ws.Range("H134").Value = 68139
$ws.Range("J134").Value = 68139
$ws.Range("L134").Value = 204417
$ws.Range("N134").Value = -209487
$ws.Range("H136").Value = 838.28125
$ws.Range("I136").Value = 634.8570999999999
$ws.Range("J136").Value = 2262.25
$ws.Range("K136").Value = 1904.5713
$ws.Range("L136").Value = 6786.75
$ws.Range("M136").Value = 645.4287000000002
$ws.Range("N136").Value = -11886.75
